$d = $word.ActiveDocument

# Use track-changes so Word's Find/Replace records the edit as a
# surgical insertion/deletion around the existing runs instead of
# collapsing (merging) the whole paragraph into a single fresh run.
# Accepting the revisions afterwards then yields the final text while
# preserving sibling (e.g. empty) runs that Find/Replace would
# otherwise prune away.
$d.TrackRevisions = $true

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Title (appears twice: main heading + bold "review" line near the end)
Replace-Text "Play Eagle Sun for Free - Exciting Features and Immersive Gameplay" "Play Eagle Sun Free - Review & Gameplay"
Replace-Text "Play Eagle Sun for Free - Exciting Features and Immersive Gameplay" "Play Eagle Sun Free - Review & Gameplay"

# "What we like" bullet list
Replace-Text "Impressive artistic production of the old west theme" "Impressive artistic production"
Replace-Text "Special features and game modes offer exciting spins" "Immersive old west theme"
Replace-Text "Immersive gameplay with appropriate music and sound effects" "Several special spins and game modes"
Replace-Text "Decent RTP provides a balanced payout experience" "Decent RTP"

# "What we don't like" bullet list
Replace-Text "No progressive jackpot feature" "Limited variety in symbols"
Replace-Text "Limited symbol variety" "Average RTP"

# Meta description (italic line)
Replace-Text "Experience the old west with Eagle Sun's impressive graphics and special features. Enjoy balanced payouts with an RTP of 95.76%. Play now for free." "Read our review of Eagle Sun and play the game for free. Enjoy an immersive old west experience with special spins and decent RTP."

# Finalize: turn tracking off and accept everything so the saved
# document contains plain (non-tracked) final text, matching the diff.
$d.TrackRevisions = $false
$d.AcceptAllRevisions()
